$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Task List")

# Fill in the Status/Issue-resolution columns for the existing last row (row 9)
$ws.Range("E9").Value = "DONE"
$ws.Range("F9").Value = "NON"

# Add a new row (10) by copying the row above (keeps borders/number formats/etc
# consistent with the rest of the table) and then overwrite with the new data.
$ws.Range("A9:G9").Copy($ws.Range("A10:G10"))

$ws.Range("A10").Value = 43145
$ws.Range("B10").Value = "AS ABOVE"
$ws.Range("C10").Value = "UPLOADED TO UBA AND TESTING"
$ws.Range("C10").WrapText = $true
$ws.Range("D10").Value = 43146
$ws.Range("D10").NumberFormat = "d-mmm"
$ws.Range("E10").Value = "UBA"

$ws.Rows.Item(10).RowHeight = 16.8

$ws.Range("E13").Select()
